$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 437. This shifts rows 437..527 down to 438..528,
# exactly matching the "old row N" -> "new row N+1" pattern seen throughout
# the diff, and leaves a brand-new, empty row 437 ready to be filled in.
$ws.Rows.Item(437).Insert()

# Populate the newly inserted row 437 with its data (columns that are
# unchanged relative to their neighbours are also (re)written here so the
# row is fully specified).
$ws.Range("A437").Value = 4
$ws.Range("B437").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C437").Value = "Los Lagos"
$ws.Range("D437").Value = 45173
$ws.Range("E437").Value = 10
$ws.Range("F437").Value = 100114014
$ws.Range("G437").Value = "Betarraga"
$ws.Range("H437").Value = "Sin especificar"
$ws.Range("I437").Value = "Primera"
$ws.Range("J437").Value = 500
$ws.Range("K437").Value = 1100
$ws.Range("L437").Value = 1100
$ws.Range("M437").Value = 1100
$ws.Range("N437").Value = '$/paquete 5 unidades'
$ws.Range("O437").Value = "Región Metropolitana"
$ws.Range("P437").Value = 220
$ws.Range("Q437").Value = 5
$ws.Range("R437").Value = "Hortaliza"
